$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclaimer string (A59) - date change 2021-05-04 -> 2021-05-05
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."

# Update weight (D) and percent change (E) values for rows 2-56
$ws.Range("D2").Value = 0.01123031859149725
$ws.Range("E2").Value = 0.1999306999307
$ws.Range("D3").Value = 0.01023436577551237
$ws.Range("E3").Value = 0.04067266327727803
$ws.Range("D4").Value = 0.0108974448773804
$ws.Range("E4").Value = -0.03256611165524004
$ws.Range("D5").Value = 0.01167014772077121
$ws.Range("E5").Value = -0.01269935026580049
$ws.Range("D6").Value = 0.01102285621656671
$ws.Range("E6").Value = 0.006535947712418277
$ws.Range("D7").Value = 0.01238114551695599
$ws.Range("E7").Value = 0.02426343154246124
$ws.Range("D8").Value = 0.01115894086121568
$ws.Range("E8").Value = 0.002450980392156854
$ws.Range("D9").Value = 0.01137151752834708
$ws.Range("E9").Value = -0.01337504888541252
$ws.Range("D10").Value = 0.01073467696907477
$ws.Range("E10").Value = -0.004349987571464098
$ws.Range("D11").Value = 0.01116694584031267
$ws.Range("E11").Value = 0.007616487455197163
$ws.Range("D12").Value = 0.4434625003418793
$ws.Range("E12").Value = 0.01875808538163004
$ws.Range("D13").Value = 0.01201325001928978
$ws.Range("E13").Value = -0.01667715544367521
$ws.Range("D14").Value = 0.01104976184075384
$ws.Range("E14").Value = -0.01803070855050004
$ws.Range("D15").Value = 0.01046384184295966
$ws.Range("E15").Value = -0.01542776998597484
$ws.Range("D16").Value = 0.01005380902476896
$ws.Range("E16").Value = -0.02932719953996543
$ws.Range("D17").Value = 0.009875587059595512
$ws.Range("E17").Value = -0.01919504643962833
$ws.Range("D18").Value = 0.008873296968492182
$ws.Range("E18").Value = 0.001729106628241883
$ws.Range("D19").Value = 0.009072087282734306
$ws.Range("E19").Value = 0.02156915610676724
$ws.Range("D20").Value = 0.01275682363318873
$ws.Range("E20").Value = -0.002091685549938882
$ws.Range("D21").Value = 0.01196744375001251
$ws.Range("E21").Value = 0.006521739130435078
$ws.Range("D22").Value = 0.01192352754524426
$ws.Range("E22").Value = -0.009958506224066466
$ws.Range("D23").Value = 0.0118540398794717
$ws.Range("E23").Value = 0.003235790658413062
$ws.Range("D24").Value = 0.01226429505819286
$ws.Range("E24").Value = -0.001359804188196923
$ws.Range("D25").Value = 0.01273681118544624
$ws.Range("E25").Value = 0.004608938547486119
$ws.Range("D26").Value = 0.01164857874931541
$ws.Range("E26").Value = -0.01437406940785702
$ws.Range("D27").Value = 0.01029507020033127
$ws.Range("E27").Value = 0.006458022851465683
$ws.Range("D28").Value = 0.01336964925517005
$ws.Range("E28").Value = 0.06613611416026344
$ws.Range("D29").Value = 0.01073334280589193
$ws.Range("E29").Value = 0.06712243629583603
$ws.Range("D30").Value = 0.006939093894179959
$ws.Range("E30").Value = 0.01041449697979591
$ws.Range("D31").Value = 0.005074267305375055
$ws.Range("E31").Value = -0.01150306748466245
$ws.Range("D32").Value = 0.009172705422772966
$ws.Range("E32").Value = -0.02696871628910469
$ws.Range("D33").Value = 0.0110268587061152
$ws.Range("E33").Value = 0.02514619883040936
$ws.Range("D34").Value = 0.01052265620326955
$ws.Range("E34").Value = -0.01287970838396102
$ws.Range("D35").Value = 0.009593078005630615
$ws.Range("E35").Value = -0.01029159519725575
$ws.Range("D36").Value = 0.009857242315831559
$ws.Range("E36").Value = -0.0118881118881119
$ws.Range("D37").Value = 0.0102800608645244
$ws.Range("E37").Value = -0.009928295642581464
$ws.Range("D38").Value = 0.01169171669222701
$ws.Range("E38").Value = -0.0005230125523014717
$ws.Range("D39").Value = 0.01390020148087666
$ws.Range("E39").Value = 0.001919631430765456
$ws.Range("D40").Value = 0.01133949761195908
$ws.Range("E40").Value = 0.007451564828613977
$ws.Range("D41").Value = 0.01220848256504435
$ws.Range("E41").Value = 0.02171062217689079
$ws.Range("D42").Value = 0.0112486633352612
$ws.Range("E42").Value = 0.002629107981220624
$ws.Range("D43").Value = 0.01159599048385874
$ws.Range("E43").Value = 0.01361470388019059
$ws.Range("D44").Value = 0.01073745647570567
$ws.Range("E44").Value = 0.01809954751131215
$ws.Range("D45").Value = 0.01202414568528291
$ws.Range("E45").Value = -0.0002311604253352195
$ws.Range("D46").Value = 0.0112905782952552
$ws.Range("E46").Value = 0.02166377816291165
$ws.Range("D47").Value = 0.01015865201488659
$ws.Range("E47").Value = 0.006665134451850241
$ws.Range("D48").Value = 0.009405850438973042
$ws.Range("E48").Value = 0.02978723404255312
$ws.Range("D49").Value = 0.00969669801283065
$ws.Range("E49").Value = 0.007911392405063333
$ws.Range("D50").Value = 0.009589964958204005
$ws.Range("E50").Value = 0.01502504173622699
$ws.Range("D51").Value = 0.00914413209460729
$ws.Range("E51").Value = -0.008474576271186307
$ws.Range("D52").Value = 0.01003068352959986
$ws.Range("E52").Value = 0.006650410108623328
$ws.Range("D53").Value = 0.00886818267629132
$ws.Range("E53").Value = -0.007246376811594124
$ws.Range("D54").Value = 0.004224294177645162
$ws.Range("E54").Value = -0.005527043031977819
$ws.Range("D55").Value = 0.004096770413419369
$ws.Range("E55").Value = 0.01085540599218415
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0.0125441357857925

$ws.Protect()
